# Gaby time sheets | ATTENDANCE UP TO DATE NOV 29
# Adds a new meeting-date column to each of the three attendance sheets
# (SPONSOR -> column H / Nov 29, TEAM -> column L / Nov 19,
#  TA -> column G / Nov 17) and marks attendance checkmarks for that date.

$wb = $excel.ActiveWorkbook
$checkmark = [char]0x2714

# ---------------------------------------------------------------------
# SPONSOR sheet: new column H = Nov 29
# ---------------------------------------------------------------------
$sponsor = $wb.Worksheets.Item("SPONSOR")

$sponsor.Range("H3").Value = "Nov"
$sponsor.Range("H4").Value = 29

$sponsor.Range("H5").Value = $checkmark
$sponsor.Range("H6").Value = $checkmark
$sponsor.Range("H7").Value = $checkmark
$sponsor.Range("H8").Value = $checkmark
$sponsor.Range("H9").Value = $checkmark
$sponsor.Range("H10").Value = $checkmark
$sponsor.Range("H11").Value = $checkmark

# column H moves out of the "unused" validation group (C5:D11 H5:N11)
# and into the "checked" validation group (E5:G11 -> E5:H11)
$sponsorRevalidate = $sponsor.Range("E5:H11")
$sponsorRevalidate.Validation.Delete()
$sponsorRevalidate.Validation.Add(3, 1, 1, "SPONSOR!checkbox", "0")

# ---------------------------------------------------------------------
# TEAM sheet: new column L = Nov 19
# ---------------------------------------------------------------------
$team = $wb.Worksheets.Item("TEAM")

$team.Range("L3").Value = "Nov"
$team.Range("L4").Value = 19

$team.Range("L5").Value = $checkmark
$team.Range("L6").Value = $checkmark
$team.Range("L7").Value = $checkmark
$team.Range("L8").Value = $checkmark
$team.Range("L9").Value = $checkmark
$team.Range("L10").Value = $checkmark
$team.Range("L11").Value = $checkmark

# ---------------------------------------------------------------------
# TA sheet: new column G = Nov 17 (Liam Nestelroad / row 10 was absent)
# ---------------------------------------------------------------------
$ta = $wb.Worksheets.Item("TA")

$ta.Range("G3").Value = "Nov"
$ta.Range("G4").Value = 17

$ta.Range("G5").Value = $checkmark
$ta.Range("G6").Value = $checkmark
$ta.Range("G7").Value = $checkmark
$ta.Range("G8").Value = $checkmark
$ta.Range("G9").Value = $checkmark
$ta.Range("G11").Value = $checkmark
# G10 intentionally left blank (row 10 missed this meeting)

# ---------------------------------------------------------------------
# Active sheet / selection bookkeeping: TA becomes the active tab,
# with G11 selected (SPONSOR and TEAM keep a secondary selection
# anchored near the newly-filled cells).
# ---------------------------------------------------------------------
$sponsor.Range("G11").Select()
$sponsor.Range("H5").Select()

$team.Range("G11").Select()
$team.Range("L5").Select()

$ta.Activate()
$ta.Range("G11").Select()
